# Updated requirements: updated requirement details to include patient page layout
# Target sheet: "Device App" (xl/worksheets/sheet5.xml)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Device App")
$ws.Activate()

# ---------------------------------------------------------------------------
# 1) Make sure every new cell we are about to touch (rows 12-22, cols B:E)
#    starts from the same "bordered + wrap text" look already used by the
#    rest of the table (style used on row 11). Doing this first means the
#    later .Value assignments just fill text into already-styled cells.
# ---------------------------------------------------------------------------
$ws.Range("B11:E11").Copy()
$ws.Range("B12:E22").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2) Row heights for the rewritten / new rows (matches the taller wrapped
#    paragraphs that were added).
# ---------------------------------------------------------------------------
$ws.Rows.Item(12).RowHeight = 60
$ws.Rows.Item(13).RowHeight = 45
$ws.Rows.Item(14).RowHeight = 30
$ws.Rows.Item(15).RowHeight = 30
$ws.Rows.Item(16).RowHeight = 120
$ws.Rows.Item(17).RowHeight = 90
$ws.Rows.Item(18).RowHeight = 30

# ---------------------------------------------------------------------------
# 3) Cell content - "Patient Page" section describing the new tabbed layout.
# ---------------------------------------------------------------------------

# Row 12 - Patient Page / Layout / tabs overview
$ws.Range("B12").Value = "Patient Page"
$ws.Range("C12").Value = "Layout"
$ws.Range("D12").Value = "Patient page should be divided in 3 tabs,
Details,
Chart,
Reports"
$ws.Range("E12").Value = ""

# Row 13 - Action bar
$ws.Range("C13").Value = "Action bar"
$ws.Range("D13").Value = "Connectivity status
Ward Name
Bed icon and Bed number"

# Row 14 - Details -> Patient Details Card
$ws.Range("C14").Value = "Details -> Patient Details Card"
$ws.Range("D14").Value = "A card with patient name and reason for admission and hospitalization date."

# Row 15 - arrow to show all details
$ws.Range("D15").Value = "An arrow to show all details (patient and medical details captured while admission)"

# Row 16 - Details -> Charts
$ws.Range("C16").Value = "Details -> Charts"
$ws.Range("D16").Value = "This needs to be discussed furher, Initial idea is to show a series of charts based on the vital parameters being monitored. The way the different charts to be displayed needs to finalized"
$ws.Range("E16").Value = "there are multiple ways to display the charts here,
1. show all charts one by one as user scrolls.
2. show one or 2 charts and give option to show more charts
3. show only one chart and provide a option to select the parameters"

# Row 17 - default chart window options / select date
$ws.Range("D17").Value = "By default the chart will show for last 24 hours, the user needs to be provided to select following options,
1. Last 24 hours
2. till now
3. Select date"
$ws.Range("E17").Value = "Select date will give user an option to select any date between admission date and current date"

# Row 18 - TBD note on vital parameter data format (highlighted)
$ws.Range("E18").Value = "The data format and configuration required for vital parameters is TBD."

# Rows 19-20 stay blank (spacer rows already carry the bordered style from step 1)

# Rows 21-22 - old "Status" / "Data" rows, renamed and pushed down
$ws.Range("C21").Value = "Patient Chart"
$ws.Range("C22").Value = "Report"

# ---------------------------------------------------------------------------
# 4) Highlight the two TBD / open-question cells in yellow.
# ---------------------------------------------------------------------------
$ws.Range("D16").Interior.Color = 65535
$ws.Range("E18").Interior.Color = 65535

# ---------------------------------------------------------------------------
# 5) Restore the view: scrolled down a bit further, selection on C13.
# ---------------------------------------------------------------------------
$ws.Range("C13").Select()
$excel.ActiveWindow.ScrollRow = 11
$excel.ActiveWindow.ScrollColumn = 1
